$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 23.15000000000018
$ws.Range("H2").Value = 0.00386416876406892
$ws.Range("I2").Value = 0.00386416876406892
$ws.Range("L2").Value = 38.9570380676392
$ws.Range("M2").Value = "[13.975690610088392, 63.938385525190014]"
$ws.Range("N2").Value = 0.002974901926535933
$ws.Range("O2").Value = 0.002974901926535933
$ws.Range("P2").Value = 1.603816069400195
$ws.Range("Q2").Value = "[0.748447499053424, 2.459184639746966]"
$ws.Range("R2").Value = 0.0004640924917651557
$ws.Range("S2").Value = 0.0004640924917651557
$ws.Range("T2").Value = 54.86119218033259
$ws.Range("U2").Value = "[39.49325926951502, 70.22912509115017]"
$ws.Range("V2").Value = 0.000000005329617636107287
$ws.Range("W2").Value = 0.000000005329617636107287
$ws.Range("X2").Value = 17.24084084084097
$ws.Range("Y2").Value = 14.0892892892894
$ws.Range("Z2").Value = 20.39239239239255

# Row 3
$ws.Range("F3").Value = 23.15000000000018
$ws.Range("H3").Value = 0.0004398141388173293
$ws.Range("I3").Value = 0.0004398141388173293
$ws.Range("L3").Value = 41.93397949221674
$ws.Range("M3").Value = "[17.676177699294755, 66.19178128513873]"
$ws.Range("N3").Value = 0.001119990625147338
$ws.Range("O3").Value = 0.001119990625147338
$ws.Range("P3").Value = 1.830237161550811
$ws.Range("Q3").Value = "[1.1635528346628883, 2.4969214884387343]"
$ws.Range("R3").Value = 0.000001551535694144235
$ws.Range("S3").Value = 0.000001551535694144235
$ws.Range("T3").Value = 52.08092657673123
$ws.Range("U3").Value = "[38.37135575379554, 65.79049739966692]"
$ws.Range("V3").Value = 0.000000001114260683721113
$ws.Range("W3").Value = 0.000000001114260683721113
$ws.Range("X3").Value = 16.40660660660673
$ws.Range("Y3").Value = 13.95025025025036
$ws.Range("Z3").Value = 18.8629629629631

# Row 4
$ws.Range("F4").Value = 23.15000000000018
$ws.Range("H4").Value = 0.0000002968814771797668
$ws.Range("I4").Value = 0.0000002968814771797668
$ws.Range("L4").Value = 59.8959893193038
$ws.Range("M4").Value = "[37.02194240077068, 82.77003623783692]"
$ws.Range("N4").Value = 0.000003675999582730682
$ws.Range("O4").Value = 0.000003675999582730682
$ws.Range("P4").Value = 1.842816111114733
$ws.Range("Q4").Value = "[1.427710775505271, 2.257921446724196]"
$ws.Range("R4").Value = 0.00000000001533373428230789
$ws.Range("S4").Value = 0.00000000001533373428230789
$ws.Range("T4").Value = 69.24793370846726
$ws.Range("U4").Value = "[56.38829765350384, 82.10756976343067]"
$ws.Range("V4").Value = 0.00000000000003841371665203042
$ws.Range("W4").Value = 0.00000000000003841371665203042
$ws.Range("X4").Value = 16.36026026026039
$ws.Range("Y4").Value = 14.83083083083095
$ws.Range("Z4").Value = 17.88968968968983

# Row 5
$ws.Range("F5").Value = 23.15000000000018
$ws.Range("H5").Value = 0.0005252303066670683
$ws.Range("I5").Value = 0.0005252303066670683
$ws.Range("L5").Value = 37.14857352464572
$ws.Range("M5").Value = "[13.569596573854554, 60.727550475436885]"
$ws.Range("N5").Value = 0.00271740545441479
$ws.Range("O5").Value = 0.00271740545441479
$ws.Range("P5").Value = 2.182447749340658
$ws.Range("Q5").Value = "[1.5660792207084269, 2.7988162779728887]"
$ws.Range("R5").Value = 0.000000006504698113474205
$ws.Range("S5").Value = 0.000000006504698113474205
$ws.Range("T5").Value = 58.44902130655093
$ws.Range("U5").Value = "[46.011103633128215, 70.88693897997365]"
$ws.Range("V5").Value = 0.000000000002830180534374449
$ws.Range("W5").Value = 0.000000000002830180534374449
$ws.Range("X5").Value = 15.10890890890902
$ws.Range("Y5").Value = 12.83793793793804
$ws.Range("Z5").Value = 17.37987987988001

# Row 6
$ws.Range("F6").Value = 23.15000000000018
$ws.Range("H6").Value = 0.00004363399727269623
$ws.Range("I6").Value = 0.00004363399727269623
$ws.Range("L6").Value = 51.13198758032942
$ws.Range("M6").Value = "[23.74959091044495, 78.51438425021388]"
$ws.Range("N6").Value = 0.0004863945430486627
$ws.Range("O6").Value = 0.0004863945430486627
$ws.Range("P6").Value = 2.597553084950119
$ws.Range("Q6").Value = "[2.0315003545735806, 3.163605815326658]"
$ws.Range("R6").Value = 0.000000000005778488798569015
$ws.Range("S6").Value = 0.000000000005778488798569015
$ws.Range("T6").Value = 62.06438141126353
$ws.Range("U6").Value = "[47.46364482231779, 76.66511800020926]"
$ws.Range("V6").Value = 0.00000000005327871477334156
$ws.Range("W6").Value = 0.00000000005327871477334156
$ws.Range("X6").Value = 13.57947947947958
$ws.Range("Y6").Value = 11.49389389389398
$ws.Range("Z6").Value = 15.66506506506518

# Row 7
$ws.Range("F7").Value = 23.15000000000018
$ws.Range("H7").Value = 0.00000004917596185816109
$ws.Range("I7").Value = 0.00000004917596185816109
$ws.Range("L7").Value = 61.47231037764974
$ws.Range("M7").Value = "[39.15314578139068, 83.7914749739088]"
$ws.Range("N7").Value = 0.000001459479263887786
$ws.Range("O7").Value = 0.000001459479263887786
$ws.Range("P7").Value = 2.823974177100735
$ws.Range("Q7").Value = "[2.4340267406191187, 3.2139216135823507]"
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 64.39050657599476
$ws.Range("U7").Value = "[51.833841367552964, 76.94717178443656]"
$ws.Range("V7").Value = 0.0000000000001871836019518014
$ws.Range("W7").Value = 0.0000000000001871836019518014
$ws.Range("X7").Value = 12.74524524524534
$ws.Range("Y7").Value = 11.3085085085086
$ws.Range("Z7").Value = 14.18198198198209

# Row 8
$ws.Range("F8").Value = 23.15000000000018
$ws.Range("H8").Value = 0.00002849298943596867
$ws.Range("I8").Value = 0.00002849298943596867
$ws.Range("L8").Value = 39.45059510095895
$ws.Range("M8").Value = "[19.669516421908227, 59.23167378000967]"
$ws.Range("N8").Value = 0.0002213472879342149
$ws.Range("O8").Value = 0.0002213472879342149
$ws.Range("P8").Value = 2.836553126664658
$ws.Range("Q8").Value = "[2.257921446724196, 3.4151848066051205]"
$ws.Range("R8").Value = 0.0000000000007736034035588091
$ws.Range("S8").Value = 0.0000000000007736034035588091
$ws.Range("T8").Value = 53.18288476391037
$ws.Range("U8").Value = "[42.139205719687865, 64.22656380813287]"
$ws.Range("V8").Value = 0.000000000001341593502957039
$ws.Range("W8").Value = 0.000000000001341593502957039
$ws.Range("X8").Value = 12.698898898899
$ws.Range("Y8").Value = 10.56696696696705
$ws.Range("Z8").Value = 14.83083083083095

# Row 9
$ws.Range("F9").Value = 22.21000000000003
$ws.Range("H9").Value = 0.0002894806640361125
$ws.Range("I9").Value = 0.0002894806640361125
$ws.Range("L9").Value = 56.38229665040861
$ws.Range("M9").Value = "[20.73616735942737, 92.02842594138984]"
$ws.Range("N9").Value = 0.002623286734116048
$ws.Range("O9").Value = 0.002623286734116048
$ws.Range("P9").Value = 2.421447791055196
$ws.Range("Q9").Value = "[1.8553950606786556, 2.987500521431736]"
$ws.Range("R9").Value = 0.00000000004453504232060368
$ws.Range("S9").Value = 0.00000000004453504232060368
$ws.Range("T9").Value = 70.48889192064389
$ws.Range("U9").Value = "[52.060068957602454, 88.91771488368533]"
$ws.Range("V9").Value = 0.0000000009333547268397524
$ws.Range("W9").Value = 0.0000000009333547268397524
$ws.Range("X9").Value = 13.65059059059061
$ws.Range("Y9").Value = 11.6496896896897
$ws.Range("Z9").Value = 15.65149149149152

# Row 10
$ws.Range("F10").Value = 22.21000000000003
$ws.Range("H10").Value = 0.0001441171492491122
$ws.Range("I10").Value = 0.0001441171492491122
$ws.Range("L10").Value = 45.8012043604716
$ws.Range("M10").Value = "[20.920695825891272, 70.68171289505193]"
$ws.Range("N10").Value = 0.0005716222122655878
$ws.Range("O10").Value = 0.0005716222122655878
$ws.Range("P10").Value = -2.930895248394081
$ws.Range("Q10").Value = "[-3.585000625718082, -2.2767898710700805]"
$ws.Range("R10").Value = 0.0000000000116924248061423
$ws.Range("S10").Value = 0.0000000000116924248061423
$ws.Range("T10").Value = 68.83258448197918
$ws.Range("U10").Value = "[54.475156710896954, 83.19001225306141]"
$ws.Range("V10").Value = 0.000000000001538769112130467
$ws.Range("W10").Value = 0.000000000001538769112130467
$ws.Range("X10").Value = 10.36022022022024
$ws.Range("Y10").Value = 8.048068068068082
$ws.Range("Z10").Value = 12.67237237237239

# Row 11
$ws.Range("F11").Value = 22.21000000000003
$ws.Range("H11").Value = 0.000004691070541240627
$ws.Range("I11").Value = 0.000004691070541240627
$ws.Range("L11").Value = 54.49315747898922
$ws.Range("M11").Value = "[27.944528246091757, 81.04178671188669]"
$ws.Range("N11").Value = 0.000153289639055032
$ws.Range("O11").Value = 0.000153289639055032
$ws.Range("P11").Value = -2.377421467581465
$ws.Range("Q11").Value = "[-2.8302636518826954, -1.924579283280234]"
$ws.Range("R11").Value = 0.0000000000000879296635503124
$ws.Range("S11").Value = 0.0000000000000879296635503124
$ws.Range("T11").Value = 57.22759354448456
$ws.Range("U11").Value = "[43.75930639729138, 70.69588069167773]"
$ws.Range("V11").Value = 0.0000000000538866729016263
$ws.Range("W11").Value = 0.0000000000538866729016263
$ws.Range("X11").Value = 8.403783783783794
$ws.Range("Y11").Value = 6.803063063063074
$ws.Range("Z11").Value = 10.00450450450452

# Row 12
$ws.Range("F12").Value = 22.21000000000003
$ws.Range("H12").Value = 0.000006660569313421405
$ws.Range("I12").Value = 0.000006660569313421405
$ws.Range("L12").Value = 46.39793145484609
$ws.Range("M12").Value = "[23.58616522282291, 69.20969768686928]"
$ws.Range("N12").Value = 0.0001724885281040311
$ws.Range("O12").Value = 0.0001724885281040311
$ws.Range("P12").Value = -2.427737265837158
$ws.Range("Q12").Value = "[-2.905737349266236, -1.9497371824080791]"
$ws.Range("R12").Value = 0.0000000000002542410726391608
$ws.Range("S12").Value = 0.0000000000002542410726391608
$ws.Range("T12").Value = 56.75724088445714
$ws.Range("U12").Value = "[45.036391201321905, 68.47809056759237]"
$ws.Range("V12").Value = 0.000000000001131095217488109
$ws.Range("W12").Value = 0.000000000001131095217488109
$ws.Range("X12").Value = 8.581641641641655
$ws.Range("Y12").Value = 6.891991991992
$ws.Range("Z12").Value = 10.27129129129131
